$d = $word.ActiveDocument

# Find the "LOB1011: Eletricidade Aplicada ..." requirement paragraph; the
# trailing blank line, the page-break spacer paragraph, and the site's
# copyright/footer paragraph that follow it are being removed from the
# page.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*LOB1011: Eletricidade Aplicada*") {
        $target = $i
        break
    }
}

# Remove the next three paragraphs (blank line, page-break blank line, and
# the "© 2020 . Contact: luizeleno@usp.br. ..." copyright paragraph).
$d.Paragraphs($target + 1).Range.Delete()
$d.Paragraphs($target + 1).Range.Delete()
$d.Paragraphs($target + 1).Range.Delete()
